$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Data" to "Summary"
$ws.Name = "Summary"

# Re-assert formatting on the pre-existing header cells so their named
# styles ("name" / "title") survive the save round-trip untouched.
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# The old row 5/6 content (Micro/SMEs/MSMEs header + Enterprises % label) moves
# down to rows 10/11 below, so clear the old location completely first.
$ws.Range("A5:D6").Clear()

# --- New content ---
# Row 8: "Source Type..." bold + underline
$ws.Range("A8").Value = "Source Type: Ministry of Finance/Central Bank"
$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").Font.Underline = $true

# Row 10: headers (bold) - Micro / SMEs / MSMEs
$ws.Range("B10").Value = "Micro"
$ws.Range("B10").Font.Bold = $true

$ws.Range("C10").Value = "SMEs"
$ws.Range("C10").Font.Bold = $true

$ws.Range("D10").Value = "MSMEs"
$ws.Range("D10").Font.Bold = $true

# Row 11: "Enterprises (% of total)" bold label + value 93.1 (stored as text) normal style
$ws.Range("A11").Value = "Enterprises (% of total)"
$ws.Range("A11").Font.Bold = $true

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "93.1"

# Row 12: source citation (italic)
$ws.Range("A12").Value = "Source: UNICONS - Central Bank of Sudan, 2005"
$ws.Range("A12").Font.Italic = $true

# Row 19: bold label
$ws.Range("A19").Value = "UNICONS - Central Bank of Sudan"
$ws.Range("A19").Font.Bold = $true

# Row 20: full citation (italic)
$ws.Range("A20").Value = "UNICONS - Central Bank of Sudan, ""SITUATIONAL ANALYSIS OF THE MICROFINANCE SECTOR IN SUDAN"", 2006, p. 8-9. Available at http://www.mfu.gov.sd/sites/default/files/microsoft_word_-_situational_analysis_of_mf_sector_in_sudan-unicons.pdf"
$ws.Range("A20").Font.Italic = $true
